$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Host "WARNING: not found: $old"
    }
}

# ---------------------------------------------------------------
# 1. Heading2 -> Heading3 for the three section headings
# ---------------------------------------------------------------
$r1 = $d.Content
$r1.Find.Execute("Background:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r1.Paragraphs.Item(1).Style = "Heading 3"

$r2 = $d.Content
$r2.Find.Execute("Study 1 (Value 10/50):", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r2.Paragraphs.Item(1).Style = "Heading 3"

$r3 = $d.Content
$r3.Find.Execute("Study 2 (Value 40/50):", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r3.Paragraphs.Item(1).Style = "Heading 3"

# ---------------------------------------------------------------
# 2. Background paragraph: heart attacks sentence
# ---------------------------------------------------------------
Replace-Text "which can increase the risk for heart attacks --- also known as myocardial infarctions (MI)" "which can increase the risk for heart attacks."

# ---------------------------------------------------------------
# 3. n-3 LC-PUFA paragraph (only first run text changes)
# ---------------------------------------------------------------
Replace-Text "we can also synthesize them from the FA alpha-linolenic acid (ALA). The n-6 LC-PUFA equivalent of the n-3 LC-PUFA is arachidonic acid (ARA), which can be used in signaling pro-inflammatory processes. As with the n-3 LC-PUFA, ARA can be obtained from the diet as well as synthesized from linoleic acid (LA). However," "we can also synthesize them from the fatty acid (FA) alpha-linolenic acid (ALA; also an n-3). Conversely, the n-6 long chain polyunsaturated fatty acid (n-6 LC-PUFA) arachidonic acid (ARA) is used in signaling pro-inflammatory processes. Similar to the n-3 LC-PUFA, ARA can be obtained from the diet as well as synthesized from linoleic acid (LA; also is an n-6). However,"

# ---------------------------------------------------------------
# 4. Study 1 paragraph 1 (BMI / CRP wording)
# ---------------------------------------------------------------
Replace-Text "participants had their body mass determined and had blood samples taken. Blood samples were analyzed for C-reactive protein (CRP), which is a marker of systemic inflammation, and were also analyzed for serum FA and D6D activity." "participants had their body mass index (BMI) determined and had blood samples taken. Blood samples were analyzed for C-reactive protein (CRP; a protein that rises in response to, as well as contributes to, inflammation), serum FA, and D6D activity."

# ---------------------------------------------------------------
# 5. "After 15 years" paragraph (first run only)
# ---------------------------------------------------------------
Replace-Text "Relative risks (RR) were calculated on tertiles of CRP and D6D activity with CAD events. RR indicate the percent in risk that greater CRP or D6D have on CAD events. A RR is" "Relative risks (RR) were calculated on tertiles of serum CRP levels and D6D activity with CAD events. The RR represents the risk as a percent that greater CRP (or D6D) have on CAD events (for example, a RR of 1.30 equals a 30% greater risk). A RR is"

# ---------------------------------------------------------------
# 6. ImageCaption (Figure 1, D6D/CRP)
# ---------------------------------------------------------------
Replace-Text "Relative risks of tertiles of D6D and CRP with CAD. Ranges above the 1.0 line (dashed horizontal line) are considered statistically significant." "Relative risks of tertiles of D6D and CRP with CAD. Error bars (RR range) that cross the 1.0 line (dashed horizontal line) are not considered statistically significant."

$rfig1 = $d.Content
$rfig1.Find.Execute("are not considered statistically significant.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$notStart = $rfig1.Start + 4
$notRange = $d.Range($notStart, $notStart + 3)
$notRange.Font.Bold = 1

# ---------------------------------------------------------------
# 7. New empty paragraph before "Study 2" heading
# ---------------------------------------------------------------
$rp = $d.Content
$rp.Find.Execute("Drawing on information from the Background and the data in Table 1 and Figure 1, discuss a potential mechanism for how dietary PUFA may influence the risk for developing CAD. (Value 6/50)", $true, $false, $false, $false, $false, $true, 1, $false, "(Value 6/50)`r", 2)

# ---------------------------------------------------------------
# 8. Study 2 paragraph 1
# ---------------------------------------------------------------
Replace-Text 'aim to reduce dietary n-6 PUFA (indicated as the "Low" group). A nearby community with similar characteristics as the intervention community was used as the control group (indicated as the "High" group). Body mass index, dietary intake, and blood samples were collected from each participant. Blood samples were used to measure CRP and serum FA.' 'aim to reduce dietary n-6 PUFA (indicated as the "Low n-6" group). A nearby community with similar characteristics as the intervention community was used as the control group (indicated as the "High n-6" group). BMI, dietary intake, and blood samples were collected from each participant. Blood samples were used to measure serum CRP and FA levels.'

# ---------------------------------------------------------------
# 9. "These group" -> "This group"
# ---------------------------------------------------------------
Replace-Text "The field of nutrigenomics has revealed several candidate genes that may influence FA metabolism. These group of alleles, called the" "The field of nutrigenomics has revealed several candidate genes that may influence FA metabolism. This group of alleles, called the"

# ---------------------------------------------------------------
# 10. "have been associated" -> "has been associated"; Low -> Low FADS
# ---------------------------------------------------------------
Replace-Text 'gene cluster, have been associated with modulation in D6D activity. Therefore, a cheek swab was taken to extract DNA to quantify the gene cluster in the participants, who were then classified as "Low" if they had <4' 'gene cluster, has been associated with modulation in D6D activity. Therefore, a cheek swab was taken to extract DNA to quantify the gene cluster in the participants, who were then classified as "Low FADS" if they had <4'

# ---------------------------------------------------------------
# 11. High -> High FADS
# ---------------------------------------------------------------
Replace-Text 'alleles present and "High" if they had >4' 'alleles present and "High FADS" if they had >4'

# ---------------------------------------------------------------
# 12. "Low" FADS alleles -> Low FADS".
# ---------------------------------------------------------------
Replace-Text '* indicates significantly different (p<0.05) from participants with "Low" FADS alleles.' '* indicates significantly different (p<0.05) from participants with "Low FADS".'

# ---------------------------------------------------------------
# 13. has has a mixed -> has a mixed
# ---------------------------------------------------------------
Replace-Text "Imagine you are clinician and a patient comes in who has has a mixed, but predominately African ancestry." "Imagine you are clinician and a patient comes in who has a mixed, but predominately African ancestry."

# ---------------------------------------------------------------
# 14. ratios -> ratio, Inuits text expansion
# ---------------------------------------------------------------
Replace-Text "Given that all of these studies were conducted in Western countries with a high n-6 to n-3 ratios, discuss how the association between higher D6D activity and CAD risk may differ in countries with a lower dietary n-6 to n-3 ratio (for example, in Inuits consuming a traditional diet)." "Given that all of these studies were conducted in Western countries with a high n-6 to n-3 ratio, discuss how the association between higher D6D activity and CAD risk may differ in countries with a lower dietary n-6 to n-3 ratio (for example, in Inuits consuming a traditional diet with a high dietary intake of fish and seafood, which has large amounts of n-3 FA)."

# ---------------------------------------------------------------
# 15. myocardial infarction -> CAD ; supplementation -> nutrient supplementation
# ---------------------------------------------------------------
Replace-Text "A recent large randomized, controlled clinical trial showed no effect of n-3 LC-PUFA supplementation on myocardial infarction (a common outcome of CAD). Comment on 1) why a clinical trial may not always be able to identify effects with supplementation in the general population and 2) why targeting only n-3 LC-PUFA may not always be effective." "A recent large randomized, controlled clinical trial showed no effect of n-3 LC-PUFA supplementation on CAD. Comment on 1) why a clinical trial may not always be able to identify effects with nutrient supplementation in the general population and 2) why targeting only n-3 LC-PUFA may not always be effective."

# ---------------------------------------------------------------
# 16. Figure 2 caption - add (low FADS)/(high FADS)
# ---------------------------------------------------------------
Replace-Text "Effect of intervention on participants with either a low or a high number of FADS alleles. LL = low n-6 (intervention) and <4 FADS alleles (low); LH = low n-6 (intervention) and >4 FADS alleles (high); HL = high n-6 (control) and <4 FADS alleles (low); HH = high n-6 (control) and >4 FADS alleles (high)." "Effect of intervention on participants with either a low or a high number of FADS alleles. LL = low n-6 (intervention) and <4 FADS alleles (low FADS); LH = low n-6 (intervention) and >4 FADS alleles (high FADS); HL = high n-6 (control) and <4 FADS alleles (low FADS); HH = high n-6 (control) and >4 FADS alleles (high FADS)."

Write-Host "All edits applied."
